$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D/E value refresh (latest crypto price + 1h volume snapshot).
# A handful of D-column cells get a leading apostrophe so Excel keeps
# them as text (matching the source data) instead of parsing them as
# numbers and dropping trailing zeros.
$ws.Range("D2").Value = "37.180.55"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.056.31"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Value = "'248.50"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").Value = "'0.667"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'57.31"
$ws.Range("E8").Value = "  -2.00%  "
$ws.Range("D9").Value = "'0.387"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D12").Value = "'16.31"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").Value = "'0.916"
$ws.Range("E13").Value = "  +13.74%  "
$ws.Range("D14").Value = "2.355.19"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "'5.79"
$ws.Range("E15").Value = "  +3.41%  "
$ws.Range("D16").Value = "2.056.11"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'18.80"
$ws.Range("E17").Value = "  +13.83%  "
$ws.Range("D18").Value = "37.216.98"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").Value = "'74.91"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "'237.97"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'2.50"
$ws.Range("E24").Value = "  +4.28%  "
$ws.Range("D25").Value = "'9.69"
$ws.Range("E25").Value = "  +4.50%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.19"
$ws.Range("E26").Value = "  -4.79%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'170.41"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("D28").Value = "'20.29"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("D30").Value = "'5.20"
$ws.Range("E30").Value = "  +9.27%  "
$ws.Range("E31").Value = "  +2.73%  "
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").Value = "'4.63"
$ws.Range("E33").Value = "  +3.37%  "
$ws.Range("D34").Value = "'0.0885"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("E38").Value = "  -1.79%  "
$ws.Range("D39").Value = "'5.31"
$ws.Range("E39").Value = "  +18.87%  "
$ws.Range("E40").Value = "  +7.60%  "
$ws.Range("D41").Value = "'0.101"
$ws.Range("E41").Value = "  -11.76%  "
$ws.Range("D42").Value = "'17.79"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").Value = "'96.81"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").Value = "1.277.16"
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").Value = "2.240.50"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").Value = "'44.40"
$ws.Range("E51").Value = "  +0.90%  "
